$d = $word.ActiveDocument

# Normal style: East Asian font DejaVu Sans -> Tahoma
$normal = $d.Styles.Item("Normal")
$normal.Font.NameFarEast = "Tahoma"

# Heading style: East Asian font DejaVu Sans -> Tahoma
$heading = $d.Styles.Item("Heading")
$heading.Font.NameFarEast = "Tahoma"

# List style: add explicit complex-script (w:cs) font DejaVu Sans
$list = $d.Styles.Item("List")
$list.Font.NameBi = "DejaVu Sans"

# Caption style: add explicit complex-script (w:cs) font DejaVu Sans
$caption = $d.Styles.Item("Caption")
$caption.Font.NameBi = "DejaVu Sans"

# Index style: add explicit complex-script (w:cs) font DejaVu Sans
$index = $d.Styles.Item("Index")
$index.Font.NameBi = "DejaVu Sans"
